$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data for "project CRUD is done with validations"
# (order matches shared-string insertion order: update, description, code commit)
$ws.Range("F7").Value = "update"
$ws.Range("G8").Value = "ici is getting update to ck(an already exiting prj)"
$ws.Range("F3").Value = "code commit"

# Update selection to match target state
$ws.Activate()
$ws.Range("F9").Select()
